$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "o"
$ws.Range("B3").Value = "o"
$ws.Range("B4").Value = "o"
$ws.Range("B5").Value = "o"
$ws.Range("B6").Value = "o"
